$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1716
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 8089
$ws.Range("I40").Value = 2979.8
$ws.Range("K40").Value = 2979.8
$ws.Range("M40").Value = -2804.8
# Row 51: A Bile Business
$ws.Range("H51").Value = 9000
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 8500
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 8500
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -9468
# Row 99: Rumor Has It
$ws.Range("H99").Value = 504
$ws.Range("I99").Value = 504
$ws.Range("K99").Value = 1512
$ws.Range("M99").Value = -14
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 309.5
$ws.Range("I118").Value = 309.5
$ws.Range("K118").Value = 928.5
$ws.Range("M118").Value = 728.5
# Row 129: Practical Command
$ws.Range("H129").Value = 942.375
$ws.Range("I129").Value = 942.375
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2827.125
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2172.875
$ws.Range("N129").ClearContents()
# Row 138: All-night Crafting
$ws.Range("H138").Value = 3031.4614
$ws.Range("I138").Value = 1203.8889
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 3611.6667
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = 1528.3333
$ws.Range("N138").Value = -22277

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2314.9678
$ws.Range("I32").Value = 2314.9678
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2314.9678
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2027.9678
$ws.Range("N32").ClearContents()
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4997.143
$ws.Range("I61").Value = 5163.3335
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 5163.3335
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -4951.3335
$ws.Range("N61").Value = -4424
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2395.4546
$ws.Range("I74").Value = 1931.2222
$ws.Range("J74").Value = 4484.5
$ws.Range("K74").Value = 1931.2222
$ws.Range("L74").Value = 4484.5
$ws.Range("M74").Value = -1057.2222
$ws.Range("N74").Value = -6232.5
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2395.4546
$ws.Range("I77").Value = 1931.2222
$ws.Range("J77").Value = 4484.5
$ws.Range("K77").Value = 9656.110999999999
$ws.Range("L77").Value = 22422.5
$ws.Range("M77").Value = -5288.110999999999
$ws.Range("N77").Value = -31158.5
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4997.143
$ws.Range("I136").Value = 5163.3335
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 15490.0005
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -12940.0005
$ws.Range("N136").Value = -17100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 2778664.8
$ws.Range("I22").Value = 3704786.2
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 3704786.2
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -3704613.2
$ws.Range("N22").Value = -646
# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 26461.334
$ws.Range("I82").Value = 11630.4
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 11630.4
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -11247.4
$ws.Range("N82").Value = -45766
# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 26461.334
$ws.Range("I85").Value = 11630.4
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 11630.4
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -10304.4
$ws.Range("N85").Value = -47652
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2840
$ws.Range("I86").Value = 3050
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 3050
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1927
$ws.Range("N86").Value = -4246
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2840
$ws.Range("I89").Value = 3050
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 15250
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -9634
$ws.Range("N89").Value = -21232
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 11129.375
$ws.Range("J134").Value = 13607.6
$ws.Range("L134").Value = 40822.8
$ws.Range("N134").Value = -45892.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand
$ws.Range("H7").Value = 79.40000000000001
$ws.Range("J7").Value = 65.666664
$ws.Range("L7").Value = 196.999992
$ws.Range("N7").Value = -420.999992
# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 554.6667
$ws.Range("I98").Value = 516.3333
$ws.Range("J98").Value = 593
$ws.Range("K98").Value = 1548.9999
$ws.Range("L98").Value = 1779
$ws.Range("M98").Value = -50.99990000000003
$ws.Range("N98").Value = -4775
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 804.7273
$ws.Range("J113").Value = 784.375
$ws.Range("L113").Value = 2353.125
$ws.Range("N113").Value = -6693.125
# Row 129: Comfort Food
$ws.Range("H129").Value = 2753.2222
$ws.Range("I129").Value = 1968.4286
$ws.Range("J129").Value = 5500
$ws.Range("K129").Value = 5905.2858
$ws.Range("L129").Value = 16500
$ws.Range("M129").Value = -905.2857999999997
$ws.Range("N129").Value = -26500
# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 1100
$ws.Range("I140").Value = 1100
$ws.Range("K140").Value = 3300
$ws.Range("M140").Value = 1880

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 35: Necklet of Champions
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2621.3333
$ws.Range("I122").Value = 2591
$ws.Range("K122").Value = 7773
$ws.Range("M122").Value = -5323
# Row 132: On Board for Lar
$ws.Range("H132").Value = 5974.5
$ws.Range("I132").Value = 5950
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 17850
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -15320
$ws.Range("N132").Value = -23057

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 6328
$ws.Range("I7").Value = 6035
$ws.Range("K7").Value = 6035
$ws.Range("M7").Value = -5923
# Row 16: Saddle Sore
$ws.Range("H16").Value = 4268
$ws.Range("I16").Value = 4268
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4268
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4098
$ws.Range("N16").ClearContents()
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1880
$ws.Range("I22").Value = 1850
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1850
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1555
$ws.Range("N22").Value = -2590
# Row 27: Fire and Hide
$ws.Range("H27").Value = 1880
$ws.Range("I27").Value = 1850
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1850
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1743
$ws.Range("N27").Value = -2214
# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 2637
$ws.Range("I68").Value = 2637
$ws.Range("K68").Value = 2637
$ws.Range("M68").Value = -1888
# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2637
$ws.Range("I71").Value = 2637
$ws.Range("K71").Value = 13185
$ws.Range("M71").Value = -9441
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 5333.3335
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 5333.3335
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -4085.3335
$ws.Range("N93").Value = -6496
# Row 94: Fitting In
$ws.Range("H94").Value = 5037500
$ws.Range("J94").Value = 5037500
$ws.Range("L94").Value = 5037500
$ws.Range("N94").Value = -5038852
# Row 126: Battered Books
$ws.Range("H126").Value = 6328
$ws.Range("I126").Value = 6035
$ws.Range("K126").Value = 18105
$ws.Range("M126").Value = -15635

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 18: Welcome to the Cotton Club
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 22: Better Shroud than Sorry
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -307
$ws.Range("N22").ClearContents()
# Row 82: Investing in the Future
$ws.Range("H82").Value = 89999.5
$ws.Range("J82").Value = 89999.5
$ws.Range("L82").Value = 89999.5
$ws.Range("N82").Value = -90765.5
# Row 85: Maids of Honor (L)
$ws.Range("H85").Value = 89999.5
$ws.Range("J85").Value = 89999.5
$ws.Range("L85").Value = 89999.5
$ws.Range("N85").Value = -92651.5
# Row 86: Felt for the Fallen
$ws.Range("H86").Value = 100000
$ws.Range("J86").Value = 100000
$ws.Range("L86").Value = 100000
$ws.Range("N86").Value = -102246
# Row 87: He Wears the Pants
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 100000
$ws.Range("L87").Value = 100000
$ws.Range("N87").Value = -102496
# Row 89: Blinded Veil of Vigilance (L)
$ws.Range("H89").Value = 100000
$ws.Range("J89").Value = 100000
$ws.Range("L89").Value = 500000
$ws.Range("N89").Value = -511232
# Row 90: Pom Hemlock (L)
$ws.Range("H90").Value = 100000
$ws.Range("J90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("N90").Value = -312480
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3374.889
$ws.Range("I132").Value = 2258.0322
$ws.Range("K132").Value = 6774.096600000001
$ws.Range("M132").Value = -4244.096600000001
